$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entering match data for N3J13 (columns EU=minutes, EV=T/R status, EW=But/goal, EX=Passe D/assist)
$ws.Range("EV2").Value = "HG"
$ws.Range("EV3").Value = "HG"
$ws.Range("EU4").Value = 90
$ws.Range("EV4").Value = "T"
$ws.Range("EU5").Value = 90
$ws.Range("EV5").Value = "T"
$ws.Range("EV6").Value = "HG"
$ws.Range("EV7").Value = "HG"
$ws.Range("EV8").Value = "HG"
$ws.Range("EU9").Value = 85
$ws.Range("EV9").Value = "T"
$ws.Range("EU10").Value = 67
$ws.Range("EV10").Value = "T"
$ws.Range("EV11").Value = "HG"
$ws.Range("EV13").Value = "HG"
$ws.Range("EU14").Value = 45
$ws.Range("EV14").Value = "T"
$ws.Range("EU15").Value = 5
$ws.Range("EV15").Value = "R"
$ws.Range("EU16").Value = 90
$ws.Range("EV16").Value = "T"
$ws.Range("EW16").Value = 1
$ws.Range("EV17").Value = "HG"
$ws.Range("EU18").Value = 60
$ws.Range("EV18").Value = "T"
$ws.Range("EV19").Value = "NR"
$ws.Range("EU20").Value = 90
$ws.Range("EV20").Value = "T"
$ws.Range("EW20").Value = 1
$ws.Range("EX21").Value = 1
$ws.Range("EU22").Value = 90
$ws.Range("EV22").Value = "T"
$ws.Range("EU24").Value = 90
$ws.Range("EV24").Value = "T"
$ws.Range("EV25").Value = "HG"
$ws.Range("EV26").Value = "HG"
$ws.Range("EU27").Value = 90
$ws.Range("EV27").Value = "T"
$ws.Range("EU28").Value = 30
$ws.Range("EV28").Value = "R"
$ws.Range("EU29").Value = 23
$ws.Range("EV29").Value = "R"
$ws.Range("EX29").Value = 1
$ws.Range("EU30").Value = 45
$ws.Range("EV30").Value = "R"

# Update view/selection to reflect new active cell (best effort)
$ws.Range("FC25").Select()
